# Final version of sensitivity analysis used for midterm.
$wb = $excel.ActiveWorkbook

# --- Criteria sheet: tweak Cost/Risk weights, move selection ---
$criteria = $wb.Worksheets.Item("Criteria")
$criteria.Range("B4").Value = 0.15
$criteria.Range("B5").Value = 0.25

# --- Trade-off Matrix sheet: update a few scores, format the score block, move selection ---
$matrix = $wb.Worksheets.Item("Trade-off Matrix")
$matrix.Range("D2").Value = 3
$matrix.Range("E2").Value = 4
$matrix.Range("D5").Value = 4
$matrix.Range("E5").Value = 2

$scoreRange = $matrix.Range("B2:E5")
$scoreRange.Font.Name = "Arial"
$scoreRange.Font.Size = 10

# --- Restore selections / active sheet to match final saved state ---
$criteria.Activate()
[void]$criteria.Range("B6").Select()

$matrix.Activate()
[void]$matrix.Range("E7").Select()
